$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Merge "Realizado na consulta " + "<left-double-quote>" into a single run.
# -----------------------------------------------------------------------
$quoteLeft  = [char]8220
$quoteRight = [char]8221

$fullText = $d.Content.Text
$search1a = "Realizado na consulta "
$search1b = $quoteLeft
$start1a = $fullText.IndexOf($search1a)
$end1a   = $start1a + $search1a.Length
$start1b = $end1a
$end1b   = $start1b + $search1b.Length

$r1b = $d.Range($start1b, $end1b)
$r1b.Delete()
$r1a = $d.Range($start1a, $end1a)
$r1a.InsertAfter($search1b)

# -----------------------------------------------------------------------
# 2) Merge "<right-double-quote>" + ", tanto com o nome do funcionário,
#    quanto com o nome do " into a single run, WITHOUT absorbing the
#    following "gerente." run (which keeps its own formatting/rsid).
#    We briefly detune "gerente."'s color so the paragraph-level run
#    coalescing that happens on edit doesn't fold it into the merge,
#    then restore its color afterward (a pure formatting change does not
#    trigger text-run coalescing).
# -----------------------------------------------------------------------
$fullText = $d.Content.Text
$idxGerente = $fullText.IndexOf("gerente.")
$rGerente = $d.Range($idxGerente, $idxGerente + 8)
$originalColor = $rGerente.Font.Color
$rGerente.Font.Color = 0

$fullText = $d.Content.Text
$start2a = $fullText.IndexOf($quoteRight, $start1a)
$end2a   = $start2a + 1
$search2b = ", tanto com o nome do funcion" + [char]225 + "rio, quanto com o nome do "
$start2b = $end2a
$end2b   = $start2b + $search2b.Length

$r2b = $d.Range($start2b, $end2b)
$r2b.Delete()
$r2a = $d.Range($start2a, $end2a)
$r2a.InsertAfter($search2b)

$fullText = $d.Content.Text
$idxGerente = $fullText.IndexOf("gerente.")
$rGerente = $d.Range($idxGerente, $idxGerente + 8)
$rGerente.Font.Color = $originalColor

# -----------------------------------------------------------------------
# 3) Remove the two screenshot images (inline pictures) that sat in their
#    own list paragraph, leaving the (now empty) paragraph behind.
# -----------------------------------------------------------------------
while ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}
